$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format first,
# otherwise Excel auto-converts the literal (e.g. "1.00" -> 1) same as typing
# it in manually would.
$textCells = @("D4", "D5", "D6", "D11", "D15", "D19", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D39", "D40", "D41", "D42", "D43", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values exactly as scraped for this run.
$ws.Range("D2").Value = '66.443.16'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '3.587.95'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '606.21'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '148.52'
$ws.Range("D7").Value = '3.587.27'
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").Value = '7.96'
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '4.195.60'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '29.67'
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").Value = '3.579.87'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '66.486.69'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = '11.09'
$ws.Range("E19").Value = '  -2.34%  '
$ws.Range("E20").Value = '  +2.57%  '
$ws.Range("D21").Value = '14.90'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '424.06'
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").Value = '0.613'
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").Value = '78.20'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  +2.94%  '
$ws.Range("D27").Value = '9.40'
$ws.Range("E27").Value = '  +3.97%  '
$ws.Range("D28").Value = '8.14'
$ws.Range("E28").Value = '  +3.76%  '
$ws.Range("D29").Value = '2.50'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Value = '3.584.07'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").Value = '0.158'
$ws.Range("E32").Value = '  +4.42%  '
$ws.Range("D33").Value = '25.03'
$ws.Range("E33").Value = '  -1.40%  '
$ws.Range("D34").Value = '1.42'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '7.74'
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").Value = '175.37'
$ws.Range("D40").Value = '0.0855'
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").Value = '5.20'
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").Value = '0.881'
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").Value = '46.12'
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E46").Value = '  +4.88%  '
$ws.Range("D47").Value = '23.98'
$ws.Range("E47").Value = '  +3.74%  '
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '1.14'
$ws.Range("E49").Value = '  -4.24%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '7.13'
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").Value = '0.940'
$ws.Range("E51").Value = '  +1.01%  '
